# This workbook was last edited on the "MT" sheet; the author then switched
# to the "Trafos" sheet and added two new columns ("Perdas Vazio kW" and
# "Perdas Totais kW ") right before the existing "Comentario" column, filling
# in loss figures for each transformer and updating a couple of the comments.

$wb = $excel.ActiveWorkbook
$wsReg    = $wb.Worksheets.Item("Reg")
$wsMT     = $wb.Worksheets.Item("MT")
$wsTrafos = $wb.Worksheets.Item("Trafos")

# --- MT sheet: selection moved off the old A4:H4 block -------------------
$wsMT.Activate()
$wsMT.Range("E1").Select()

# --- Trafos sheet: becomes the active/selected sheet ----------------------
$wsTrafos.Activate()

# Shift the existing "Comentario" column (P) two slots to the right, to R,
# to make room for the two new columns.
$wsTrafos.Range("P1:P13").Cut($wsTrafos.Range("R1:R13"))

# The two new header cells pick up the same bold / right-bordered /
# center-top-aligned look already used for this kind of header run on the
# "Reg" sheet.
$wsReg.Range("P1").Copy()
$wsTrafos.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTrafos.Range("P1").Value = "Perdas Vazio kW"
$wsTrafos.Range("Q1").Value = "Perdas Totais kW "

# Row 2 - Substation transformer: no-load / total losses not modeled.
$wsTrafos.Range("P2").Value = 0
$wsTrafos.Range("Q2").Value = 0

# Rows 3-7 - three-phase XFM01-05 500 kVA units.
$wsTrafos.Range("P3").Value = 1.8
$wsTrafos.Range("Q3").Value = 9
$wsTrafos.Range("R3").Value = "XLT = 1 para perdas: https://www.weg.net/catalog/weg/BR/pt/IP00/Transformador-Seco-500-0kVA-13-8-0-22kV-CST-IP-00-AN/p/14543073"

$wsTrafos.Range("P4").Value = 1.8
$wsTrafos.Range("Q4").Value = 9

$wsTrafos.Range("P5").Value = 1.8
$wsTrafos.Range("Q5").Value = 9

$wsTrafos.Range("P6").Value = 1.8
$wsTrafos.Range("Q6").Value = 9

$wsTrafos.Range("P7").Value = 1.8
$wsTrafos.Range("Q7").Value = 9

# Rows 8-13 - single-phase XFM06-11 150 kVA units.
$wsTrafos.Range("P8").Value = 0.825
$wsTrafos.Range("Q8").Value = 1.28
$wsTrafos.Range("R8").Value = "XLT = 1 / Trafo monofasico, ligacao 0 porque não é importante. Perdas foram calculadas assumindo valores tipicos"

$wsTrafos.Range("P9").Value = 0.825
$wsTrafos.Range("Q9").Value = 1.28

$wsTrafos.Range("P10").Value = 0.825
$wsTrafos.Range("Q10").Value = 1.28

$wsTrafos.Range("P11").Value = 0.825
$wsTrafos.Range("Q11").Value = 1.28

$wsTrafos.Range("P12").Value = 0.825
$wsTrafos.Range("Q12").Value = 1.28

$wsTrafos.Range("P13").Value = 0.825
$wsTrafos.Range("Q13").Value = 1.28

$wsTrafos.Range("O35").Select()
